$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Remove the now-unused Sheet2 / Sheet3 tabs --------------------------
$wb.Worksheets.Item("Sheet2").Delete()
$wb.Worksheets.Item("Sheet3").Delete()

$ws = $wb.Worksheets.Item("Sheet1")

# --- Extend the bordered table formatting down to the new rows -----------
# Row 9 (A9:E9) already carries the thin-border / left-aligned style used
# by the whole table; clone it onto the newly-needed rows 11-15 before we
# fill in values so the cell styles match exactly (same style index, no
# new style entries minted).
$ws.Range("A9:E9").Copy()
$ws.Range("A11:E15").PasteSpecial(-4122)

# F14 needs the same (border-less, left aligned) numeric style as F2 / the
# old F10 - clone it across, then blank out the old F10 (value moved to F14).
$ws.Range("F2").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("F10").Clear()

# --- Row 9: temper_value -> Linkage_flag ----------------------------------
$ws.Range("B9").Value = "Linkage_flag"
$ws.Range("C9").Value = "u8"
$ws.Range("D9").Value = "联动开关"
$ws.Range("E9").Value = "0x2F08"

# --- Row 10: resetbtcnt -> SWITCHflag2 ------------------------------------
$ws.Range("B10").Value = "SWITCHflag2"
$ws.Range("C10").Value = "u8"
$ws.Range("D10").Value = "开关灯"
$ws.Range("E10").Value = "0x2F09"

# --- Row 11: all_day_micro_light_enable -----------------------------------
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "all_day_micro_light_enable"
$ws.Range("C11").Value = "u8"
$ws.Range("D11").Value = "全天伴亮开关"
$ws.Range("E11").Value = "0x2F0A"

# --- Row 12: temper_value (re-added with new address) ---------------------
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "temper_value"
$ws.Range("C12").Value = "u8"
$ws.Range("D12").Value = "冷暖度值0~100"
$ws.Range("E12").Value = "0x2F0B"

# --- Row 13: bt_and_sigmesh_duty ------------------------------------------
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "bt_and_sigmesh_duty"
$ws.Range("C13").Value = "u16"
$ws.Range("D13").Value = "蓝牙及mesh通信周期"
$ws.Range("E13").Value = "0x2F0C~0x2F0D"

# --- Row 14: resetbtcnt (moved down, value kept) --------------------------
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "resetbtcnt"
$ws.Range("C14").Value = "u8"
$ws.Range("D14").Value = "蓝牙重新连接次数"
$ws.Range("E14").Value = "0x2F80"
$ws.Range("F14").Value = 12160

# --- Row 15: bt_join_cnt ---------------------------------------------------
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "bt_join_cnt"
$ws.Range("C15").Value = "u8"
$ws.Range("D15").Value = "蓝牙配网标志"
$ws.Range("E15").Value = "0x2F81"

# --- Cosmetic touch-ups to match the author's final view ------------------
$ws.Columns.Item(2).ColumnWidth = 27.265625
$ws.Range("B22").Select()
